# Updates the worksheet date and every arithmetic problem in the table to the
# new day's values. Each problem text (e.g. "28+28=") is unique in the
# document, so a simple Find/Replace (wdReplaceOne) targeted at the whole
# document content reliably updates the correct cell.
#
# Find.Execute signature used below:
#   FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace
# (Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceOne)

$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-12 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-10-13 Monday", 2) | Out-Null
$d.Content.Find.Execute("28+28=", $true, $true, $false, $false, $false, $true, 1, $false, "86-47=", 2) | Out-Null
$d.Content.Find.Execute("8+25=", $true, $true, $false, $false, $false, $true, 1, $false, "32-16=", 2) | Out-Null
$d.Content.Find.Execute("14+77=", $true, $true, $false, $false, $false, $true, 1, $false, "12+74=", 2) | Out-Null
$d.Content.Find.Execute("83-33=", $true, $true, $false, $false, $false, $true, 1, $false, "27+20=", 2) | Out-Null
$d.Content.Find.Execute("12+30=", $true, $true, $false, $false, $false, $true, 1, $false, "3+90=", 2) | Out-Null
$d.Content.Find.Execute("42-42=", $true, $true, $false, $false, $false, $true, 1, $false, "49+9=", 2) | Out-Null
$d.Content.Find.Execute("78+9=", $true, $true, $false, $false, $false, $true, 1, $false, "86-27=", 2) | Out-Null
$d.Content.Find.Execute("57+28=", $true, $true, $false, $false, $false, $true, 1, $false, "27+49=", 2) | Out-Null
$d.Content.Find.Execute("2+32=", $true, $true, $false, $false, $false, $true, 1, $false, "71-2=", 2) | Out-Null
$d.Content.Find.Execute("38-3=", $true, $true, $false, $false, $false, $true, 1, $false, "35+22=", 2) | Out-Null
$d.Content.Find.Execute("30+64=", $true, $true, $false, $false, $false, $true, 1, $false, "91-69=", 2) | Out-Null
$d.Content.Find.Execute("13+58=", $true, $true, $false, $false, $false, $true, 1, $false, "38-28=", 2) | Out-Null
$d.Content.Find.Execute("55+29=", $true, $true, $false, $false, $false, $true, 1, $false, "86-80=", 2) | Out-Null
$d.Content.Find.Execute("90-19=", $true, $true, $false, $false, $false, $true, 1, $false, "12+64=", 2) | Out-Null
$d.Content.Find.Execute("71+3=", $true, $true, $false, $false, $false, $true, 1, $false, "56+8=", 2) | Out-Null
$d.Content.Find.Execute("18-10=", $true, $true, $false, $false, $false, $true, 1, $false, "65+10=", 2) | Out-Null
$d.Content.Find.Execute("87-33=", $true, $true, $false, $false, $false, $true, 1, $false, "65+17=", 2) | Out-Null
$d.Content.Find.Execute("2+63=", $true, $true, $false, $false, $false, $true, 1, $false, "83+11=", 2) | Out-Null
$d.Content.Find.Execute("43+26=", $true, $true, $false, $false, $false, $true, 1, $false, "34+55=", 2) | Out-Null
$d.Content.Find.Execute("32+43=", $true, $true, $false, $false, $false, $true, 1, $false, "16+28=", 2) | Out-Null
$d.Content.Find.Execute("2+56=", $true, $true, $false, $false, $false, $true, 1, $false, "98+1=", 2) | Out-Null
$d.Content.Find.Execute("25-6=", $true, $true, $false, $false, $false, $true, 1, $false, "5+58=", 2) | Out-Null
$d.Content.Find.Execute("44-27=", $true, $true, $false, $false, $false, $true, 1, $false, "61-34=", 2) | Out-Null
$d.Content.Find.Execute("13-11=", $true, $true, $false, $false, $false, $true, 1, $false, "32+56=", 2) | Out-Null
$d.Content.Find.Execute("88-41=", $true, $true, $false, $false, $false, $true, 1, $false, "3+29=", 2) | Out-Null
$d.Content.Find.Execute("22-8=", $true, $true, $false, $false, $false, $true, 1, $false, "73-41=", 2) | Out-Null
$d.Content.Find.Execute("58+37=", $true, $true, $false, $false, $false, $true, 1, $false, "78-56=", 2) | Out-Null
$d.Content.Find.Execute("67-31=", $true, $true, $false, $false, $false, $true, 1, $false, "18-4=", 2) | Out-Null
$d.Content.Find.Execute("12+61=", $true, $true, $false, $false, $false, $true, 1, $false, "37+44=", 2) | Out-Null
$d.Content.Find.Execute("67+19=", $true, $true, $false, $false, $false, $true, 1, $false, "73+12=", 2) | Out-Null
$d.Content.Find.Execute("75-71=", $true, $true, $false, $false, $false, $true, 1, $false, "32-3=", 2) | Out-Null
$d.Content.Find.Execute("69-2=", $true, $true, $false, $false, $false, $true, 1, $false, "86-22=", 2) | Out-Null
$d.Content.Find.Execute("42-2=", $true, $true, $false, $false, $false, $true, 1, $false, "82-6=", 2) | Out-Null
$d.Content.Find.Execute("49+10=", $true, $true, $false, $false, $false, $true, 1, $false, "64-33=", 2) | Out-Null
$d.Content.Find.Execute("87-65=", $true, $true, $false, $false, $false, $true, 1, $false, "97-12=", 2) | Out-Null
$d.Content.Find.Execute("18-17=", $true, $true, $false, $false, $false, $true, 1, $false, "99-19=", 2) | Out-Null
$d.Content.Find.Execute("94-62=", $true, $true, $false, $false, $false, $true, 1, $false, "87-82=", 2) | Out-Null
$d.Content.Find.Execute("16+27=", $true, $true, $false, $false, $false, $true, 1, $false, "85-27=", 2) | Out-Null
$d.Content.Find.Execute("30+4=", $true, $true, $false, $false, $false, $true, 1, $false, "52+40=", 2) | Out-Null
$d.Content.Find.Execute("99-78=", $true, $true, $false, $false, $false, $true, 1, $false, "90-44=", 2) | Out-Null
$d.Content.Find.Execute("76-57=", $true, $true, $false, $false, $false, $true, 1, $false, "25+1=", 2) | Out-Null
$d.Content.Find.Execute("23+64=", $true, $true, $false, $false, $false, $true, 1, $false, "28+24=", 2) | Out-Null
$d.Content.Find.Execute("75+11=", $true, $true, $false, $false, $false, $true, 1, $false, "81-43=", 2) | Out-Null
$d.Content.Find.Execute("8+90=", $true, $true, $false, $false, $false, $true, 1, $false, "70-55=", 2) | Out-Null
$d.Content.Find.Execute("43+39=", $true, $true, $false, $false, $false, $true, 1, $false, "46-16=", 2) | Out-Null
$d.Content.Find.Execute("37+19=", $true, $true, $false, $false, $false, $true, 1, $false, "99-95=", 2) | Out-Null
$d.Content.Find.Execute("29-16=", $true, $true, $false, $false, $false, $true, 1, $false, "98-54=", 2) | Out-Null
$d.Content.Find.Execute("11+61=", $true, $true, $false, $false, $false, $true, 1, $false, "91-58=", 2) | Out-Null
$d.Content.Find.Execute("56-43=", $true, $true, $false, $false, $false, $true, 1, $false, "41+45=", 2) | Out-Null
$d.Content.Find.Execute("16-12=", $true, $true, $false, $false, $false, $true, 1, $false, "71+24=", 2) | Out-Null
$d.Content.Find.Execute("3+8=", $true, $true, $false, $false, $false, $true, 1, $false, "43+17=", 2) | Out-Null
$d.Content.Find.Execute("24+19=", $true, $true, $false, $false, $false, $true, 1, $false, "8+58=", 2) | Out-Null
$d.Content.Find.Execute("15+69=", $true, $true, $false, $false, $false, $true, 1, $false, "31+31=", 2) | Out-Null
$d.Content.Find.Execute("83-22=", $true, $true, $false, $false, $false, $true, 1, $false, "25+3=", 2) | Out-Null
$d.Content.Find.Execute("48-35=", $true, $true, $false, $false, $false, $true, 1, $false, "89-56=", 2) | Out-Null
$d.Content.Find.Execute("51-46=", $true, $true, $false, $false, $false, $true, 1, $false, "25+15=", 2) | Out-Null
$d.Content.Find.Execute("54-15=", $true, $true, $false, $false, $false, $true, 1, $false, "36+25=", 2) | Out-Null
$d.Content.Find.Execute("49-29=", $true, $true, $false, $false, $false, $true, 1, $false, "12+39=", 2) | Out-Null
$d.Content.Find.Execute("44+18=", $true, $true, $false, $false, $false, $true, 1, $false, "62+4=", 2) | Out-Null
$d.Content.Find.Execute("91-31=", $true, $true, $false, $false, $false, $true, 1, $false, "48-25=", 2) | Out-Null
$d.Content.Find.Execute("86+12=", $true, $true, $false, $false, $false, $true, 1, $false, "24-1=", 2) | Out-Null
$d.Content.Find.Execute("20-10=", $true, $true, $false, $false, $false, $true, 1, $false, "36-30=", 2) | Out-Null
$d.Content.Find.Execute("28-15=", $true, $true, $false, $false, $false, $true, 1, $false, "63+5=", 2) | Out-Null
$d.Content.Find.Execute("45+45=", $true, $true, $false, $false, $false, $true, 1, $false, "71-15=", 2) | Out-Null
$d.Content.Find.Execute("27+21=", $true, $true, $false, $false, $false, $true, 1, $false, "58-2=", 2) | Out-Null
$d.Content.Find.Execute("27+30=", $true, $true, $false, $false, $false, $true, 1, $false, "78+0=", 2) | Out-Null
$d.Content.Find.Execute("16-0=", $true, $true, $false, $false, $false, $true, 1, $false, "63-41=", 2) | Out-Null
$d.Content.Find.Execute("55+34=", $true, $true, $false, $false, $false, $true, 1, $false, "89-42=", 2) | Out-Null
$d.Content.Find.Execute("4-2=", $true, $true, $false, $false, $false, $true, 1, $false, "41-19=", 2) | Out-Null
$d.Content.Find.Execute("97-34=", $true, $true, $false, $false, $false, $true, 1, $false, "95-7=", 2) | Out-Null
$d.Content.Find.Execute("12+72=", $true, $true, $false, $false, $false, $true, 1, $false, "51+27=", 2) | Out-Null
$d.Content.Find.Execute("34+57=", $true, $true, $false, $false, $false, $true, 1, $false, "72-15=", 2) | Out-Null
$d.Content.Find.Execute("87-4=", $true, $true, $false, $false, $false, $true, 1, $false, "20+3=", 2) | Out-Null
$d.Content.Find.Execute("45+25=", $true, $true, $false, $false, $false, $true, 1, $false, "23+58=", 2) | Out-Null
$d.Content.Find.Execute("23-15=", $true, $true, $false, $false, $false, $true, 1, $false, "12+7=", 2) | Out-Null
$d.Content.Find.Execute("61-55=", $true, $true, $false, $false, $false, $true, 1, $false, "44+33=", 2) | Out-Null
$d.Content.Find.Execute("36+59=", $true, $true, $false, $false, $false, $true, 1, $false, "63-17=", 2) | Out-Null
$d.Content.Find.Execute("18+53=", $true, $true, $false, $false, $false, $true, 1, $false, "53+44=", 2) | Out-Null
$d.Content.Find.Execute("43-3=", $true, $true, $false, $false, $false, $true, 1, $false, "83-25=", 2) | Out-Null
$d.Content.Find.Execute("51-47=", $true, $true, $false, $false, $false, $true, 1, $false, "65+1=", 2) | Out-Null
$d.Content.Find.Execute("29+57=", $true, $true, $false, $false, $false, $true, 1, $false, "44-37=", 2) | Out-Null
$d.Content.Find.Execute("83-72=", $true, $true, $false, $false, $false, $true, 1, $false, "42-29=", 2) | Out-Null
$d.Content.Find.Execute("18+78=", $true, $true, $false, $false, $false, $true, 1, $false, "8+68=", 2) | Out-Null
$d.Content.Find.Execute("2+7=", $true, $true, $false, $false, $false, $true, 1, $false, "36+62=", 2) | Out-Null
$d.Content.Find.Execute("3+63=", $true, $true, $false, $false, $false, $true, 1, $false, "8+53=", 2) | Out-Null
$d.Content.Find.Execute("34+50=", $true, $true, $false, $false, $false, $true, 1, $false, "16+20=", 2) | Out-Null
$d.Content.Find.Execute("37+46=", $true, $true, $false, $false, $false, $true, 1, $false, "0+37=", 2) | Out-Null
$d.Content.Find.Execute("80+5=", $true, $true, $false, $false, $false, $true, 1, $false, "92-84=", 2) | Out-Null
$d.Content.Find.Execute("44+31=", $true, $true, $false, $false, $false, $true, 1, $false, "40-13=", 2) | Out-Null
$d.Content.Find.Execute("18+14=", $true, $true, $false, $false, $false, $true, 1, $false, "96-0=", 2) | Out-Null
$d.Content.Find.Execute("84-61=", $true, $true, $false, $false, $false, $true, 1, $false, "80-49=", 2) | Out-Null
$d.Content.Find.Execute("73-39=", $true, $true, $false, $false, $false, $true, 1, $false, "39+10=", 2) | Out-Null
$d.Content.Find.Execute("34+44=", $true, $true, $false, $false, $false, $true, 1, $false, "50-19=", 2) | Out-Null
$d.Content.Find.Execute("81-0=", $true, $true, $false, $false, $false, $true, 1, $false, "99-38=", 2) | Out-Null
$d.Content.Find.Execute("63+1=", $true, $true, $false, $false, $false, $true, 1, $false, "11+71=", 2) | Out-Null
$d.Content.Find.Execute("80-59=", $true, $true, $false, $false, $false, $true, 1, $false, "19+12=", 2) | Out-Null
$d.Content.Find.Execute("34+62=", $true, $true, $false, $false, $false, $true, 1, $false, "24+39=", 2) | Out-Null
$d.Content.Find.Execute("4+84=", $true, $true, $false, $false, $false, $true, 1, $false, "98-30=", 2) | Out-Null
$d.Content.Find.Execute("59-21=", $true, $true, $false, $false, $false, $true, 1, $false, "55+3=", 2) | Out-Null
$d.Content.Find.Execute("39+29=", $true, $true, $false, $false, $false, $true, 1, $false, "18+38=", 2) | Out-Null
